# Generate Report for Handback
#
# The localization "Handback" pass failed to transform the handback file for
# e2e\28e34605-f3b1-4f4a-80a2-b98cc8c0d1fe.md (row 3) in both the zh-cn and
# de-de targets. Update the status everywhere it is surfaced, and record the
# error detail message for each locale.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# --- Overview sheet: zh-cn (E3) and de-de (F3) status columns for the
#     28e34605-... row -----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status (C3) + Error Detail (P3) for the same row --------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("P3").Value = "Handback file name: rxgxnuqp.mhy is different with handoff file name: 28e34605-f3b1-4f4a-80a2-b98cc8c0d1fe.3b7f24fb7bf73b94877b15d809a76be260cefca5.zh-cn."

# Widen the Error Detail column so the new message is readable (matches the
# target column width of 40 "characters"; Excel stores width + 5/6 internally
# so we back that offset out here).
$wsZh.Range("P1").ColumnWidth = 40 - (5/6)

# --- de-de sheet: Status (C3) + Error Detail (P3) for the same row --------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("P3").Value = "Handback file name: rxgxnuqp.mhy is different with handoff file name: 28e34605-f3b1-4f4a-80a2-b98cc8c0d1fe.3b7f24fb7bf73b94877b15d809a76be260cefca5.de-de."

$wsDe.Range("P1").ColumnWidth = 40 - (5/6)
